# Auto-generated Excel COM-interop script applying the diff to Carbuncle_Profits
# (workbook sheets: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 5436
$ws.Range("I19").Value = 11792.333
$ws.Range("J19").Value = 1035.4615
$ws.Range("K19").Value = 11792.333
$ws.Range("L19").Value = 1035.4615
$ws.Range("M19").Value = -11617.333
$ws.Range("N19").Value = -1385.4615
$ws.Range("H32").Value = 1730.125
$ws.Range("I32").Value = 1101
$ws.Range("J32").Value = 1820
$ws.Range("K32").Value = 1101
$ws.Range("L32").Value = 1820
$ws.Range("M32").Value = -775
$ws.Range("N32").Value = -2472
$ws.Range("H43").Value = 1712
$ws.Range("I43").Value = 1390
$ws.Range("J43").Value = 3000
$ws.Range("K43").Value = 1390
$ws.Range("L43").Value = 3000
$ws.Range("M43").Value = -1321
$ws.Range("N43").Value = -3138
$ws.Range("H88").Value = 1564
$ws.Range("J88").Value = 1591.3334
$ws.Range("L88").Value = 1591.3334
$ws.Range("N88").Value = -2403.3334
$ws.Range("H91").Value = 1564
$ws.Range("J91").Value = 1591.3334
$ws.Range("L91").Value = 1591.3334
$ws.Range("N91").Value = -4399.3334
$ws.Range("H94").Value = 3500
$ws.Range("I94").Value = 3500
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 3500
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -3049
$ws.Range("N94").ClearContents()
$ws.Range("H116").Value = 2132.5
$ws.Range("I116").Value = 1946
$ws.Range("J116").Value = 2598.75
$ws.Range("K116").Value = 1946
$ws.Range("L116").Value = 2598.75
$ws.Range("M116").Value = 1496
$ws.Range("N116").Value = -9482.75
$ws.Range("H129").Value = 867.5606
$ws.Range("J129").Value = 880.7258
$ws.Range("L129").Value = 2642.1774
$ws.Range("N129").Value = -12642.1774
$ws.Range("H132").Value = 149429.14
$ws.Range("I132").Value = 260751
$ws.Range("K132").Value = 782253
$ws.Range("M132").Value = -779723
$ws.Range("H141").Value = 10560.125
$ws.Range("I141").Value = 1669.625
$ws.Range("J141").Value = 19450.625
$ws.Range("K141").Value = 5008.875
$ws.Range("L141").Value = 58351.875
$ws.Range("M141").Value = 171.125
$ws.Range("N141").Value = -68711.875
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7831.2974
$ws.Range("I32").Value = 5517.3076
$ws.Range("J32").Value = 13300.728
$ws.Range("K32").Value = 5517.3076
$ws.Range("L32").Value = 13300.728
$ws.Range("M32").Value = -5230.3076
$ws.Range("N32").Value = -13874.728
$ws.Range("H88").Value = 10573.5
$ws.Range("I88").Value = 12188.2
$ws.Range("J88").Value = 2500
$ws.Range("K88").Value = 12188.2
$ws.Range("L88").Value = 2500
$ws.Range("M88").Value = -11782.2
$ws.Range("N88").Value = -3312
$ws.Range("H91").Value = 10573.5
$ws.Range("I91").Value = 12188.2
$ws.Range("J91").Value = 2500
$ws.Range("K91").Value = 12188.2
$ws.Range("L91").Value = 2500
$ws.Range("M91").Value = -10784.2
$ws.Range("N91").Value = -5308
$ws.Range("H102").Value = 1610.7826
$ws.Range("I102").Value = 1486.7368
$ws.Range("J102").Value = 2200
$ws.Range("K102").Value = 1486.7368
$ws.Range("L102").Value = 2200
$ws.Range("M102").Value = 135.2632000000001
$ws.Range("N102").Value = -5444
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3142.2
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 3142.2
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("H99").Value = 1363
$ws.Range("I99").Value = 899.3077
$ws.Range("J99").Value = 1911
$ws.Range("K99").Value = 899.3077
$ws.Range("L99").Value = 1911
$ws.Range("M99").Value = 598.6923
$ws.Range("N99").Value = -4907
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H36").Value = 15725
$ws.Range("I36").Value = 1500
$ws.Range("J36").Value = 29950
$ws.Range("K36").Value = 1500
$ws.Range("L36").Value = 29950
$ws.Range("M36").Value = -1112
$ws.Range("N36").Value = -30726
$ws.Range("H40").Value = 15725
$ws.Range("I40").Value = 1500
$ws.Range("J40").Value = 29950
$ws.Range("K40").Value = 1500
$ws.Range("L40").Value = 29950
$ws.Range("M40").Value = -1340
$ws.Range("N40").Value = -30270
$ws.Range("H99").Value = 7411.4116
$ws.Range("I99").Value = 1545.6923
$ws.Range("J99").Value = 26475
$ws.Range("K99").Value = 1545.6923
$ws.Range("L99").Value = 26475
$ws.Range("M99").Value = -47.69229999999993
$ws.Range("N99").Value = -29471
$ws.Range("H126").Value = 7411.4116
$ws.Range("I126").Value = 1545.6923
$ws.Range("J126").Value = 26475
$ws.Range("K126").Value = 4637.0769
$ws.Range("L126").Value = 79425
$ws.Range("M126").Value = -2167.0769
$ws.Range("N126").Value = -84365
$ws.Range("H132").Value = 3045.4443
$ws.Range("I132").Value = 1788
$ws.Range("K132").Value = 5364
$ws.Range("M132").Value = -2834
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1401.375
$ws.Range("I102").Value = 1242.2
$ws.Range("J102").Value = 1666.6666
$ws.Range("K102").Value = 1242.2
$ws.Range("L102").Value = 1666.6666
$ws.Range("M102").Value = 379.8
$ws.Range("N102").Value = -4910.6666
$ws.Range("H126").Value = 1922.7222
$ws.Range("I126").Value = 1567.7916
$ws.Range("J126").Value = 2632.5833
$ws.Range("K126").Value = 4703.3748
$ws.Range("L126").Value = 7897.749899999999
$ws.Range("M126").Value = -2233.3748
$ws.Range("N126").Value = -12837.7499
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1886.8
$ws.Range("J16").Value = 2410.2
$ws.Range("L16").Value = 2410.2
$ws.Range("N16").Value = -2750.2
$ws.Range("H46").Value = 113677.78
$ws.Range("I46").Value = 202620.2
$ws.Range("J46").Value = 2499.75
$ws.Range("K46").Value = 202620.2
$ws.Range("L46").Value = 2499.75
$ws.Range("M46").Value = -202432.2
$ws.Range("N46").Value = -2875.75
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 40060.73
$ws.Range("I122").Value = 45173.652
$ws.Range("J122").Value = 861.6667
$ws.Range("K122").Value = 135520.956
$ws.Range("L122").Value = 2585.0001
$ws.Range("M122").Value = -133070.956
$ws.Range("N122").Value = -7485.0001
$ws.Range("H126").Value = 94946
$ws.Range("I126").Value = 146385.86
$ws.Range("J126").Value = 4926.25
$ws.Range("K126").Value = 439157.58
$ws.Range("L126").Value = 14778.75
$ws.Range("M126").Value = -436687.58
$ws.Range("N126").Value = -19718.75
